$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Append a new titration/accuracy data row (row 56) below the existing data.
$ws.Range("A56").Value = 20220203
$ws.Range("B56").Value = 2226.15763917096
$ws.Range("C56").Value = 2224.4699999999998
$ws.Range("D56").Formula = "=100*(B56-C56)/C56"
$ws.Range("E56").Value = 180
$ws.Range("F56").Value = "CRM OPENED 20220118"

# Keep the selection / view in sync with the newly added row.
$ws.Range("G56").Select()
